# Refresh the cryptos price table (column D = Price, column E =
# Volume(1h)) with the latest scrape. The sheet stores every value as
# text (even price figures that look numeric, e.g. "228.63" or
# "0.0840"), so plain numeric-looking replacements must be written
# back as Text to keep exact formatting (trailing zeros, the
# thousands-dot grouping like "2.403.52", etc.) instead of being
# auto-converted to a Double by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

# Cells whose new text is not number-like (contains letters, %, spaces,
# multiple dots, or the subscript digit) -- a direct .Value assignment
# already stores these as Text.
$textUpdates = @(
    @{ Row = 2; Col = 4; Text = "38.897.81" },
    @{ Row = 2; Col = 5; Text = "  +3.03%  " },
    @{ Row = 3; Col = 4; Text = "2.093.68" },
    @{ Row = 3; Col = 5; Text = "  +2.18%  " },
    @{ Row = 4; Col = 5; Text = "  -0.04%  " },
    @{ Row = 5; Col = 5; Text = "  +0.47%  " },
    @{ Row = 6; Col = 5; Text = "  +0.81%  " },
    @{ Row = 7; Col = 5; Text = "  +0.63%  " },
    @{ Row = 9; Col = 5; Text = "  +2.28%  " },
    @{ Row = 10; Col = 5; Text = "  +0.50%  " },
    @{ Row = 11; Col = 5; Text = "  -0.07%  " },
    @{ Row = 12; Col = 4; Text = "2.403.52" },
    @{ Row = 12; Col = 5; Text = "  +2.23%  " },
    @{ Row = 13; Col = 5; Text = "  +3.94%  " },
    @{ Row = 14; Col = 5; Text = "  +2.70%  " },
    @{ Row = 15; Col = 5; Text = "  +4.36%  " },
    @{ Row = 16; Col = 5; Text = "  -0.70%  " },
    @{ Row = 17; Col = 4; Text = "2.097.13" },
    @{ Row = 17; Col = 5; Text = "  +2.57%  " },
    @{ Row = 18; Col = 4; Text = "38.806.05" },
    @{ Row = 19; Col = 5; Text = "  +3.14%  " },
    @{ Row = 20; Col = 5; Text = "  +2.15%  " },
    @{ Row = 21; Col = 4; Text = ("0.0{0}0838" -f $sub3) },
    @{ Row = 21; Col = 5; Text = "  +1.26%  " },
    @{ Row = 22; Col = 5; Text = "  +2.35%  " },
    @{ Row = 23; Col = 5; Text = "  -0.43%  " },
    @{ Row = 24; Col = 5; Text = "  -0.37%  " },
    @{ Row = 25; Col = 5; Text = "  +2.78%  " },
    @{ Row = 26; Col = 5; Text = "  +1.15%  " },
    @{ Row = 27; Col = 5; Text = "  +2.05%  " },
    @{ Row = 28; Col = 5; Text = "  +8.94%  " },
    @{ Row = 29; Col = 5; Text = "  +14.61%  " },
    @{ Row = 30; Col = 5; Text = "  +2.34%  " },
    @{ Row = 31; Col = 5; Text = "  +0.95%  " },
    @{ Row = 32; Col = 5; Text = "  +5.52%  " },
    @{ Row = 33; Col = 5; Text = "  +2.70%  " },
    @{ Row = 34; Col = 5; Text = "  +4.17%  " },
    @{ Row = 35; Col = 5; Text = "  +1.35%  " },
    @{ Row = 36; Col = 5; Text = "  -0.22%  " },
    @{ Row = 37; Col = 5; Text = "  +1.64%  " },
    @{ Row = 38; Col = 5; Text = "  +3.80%  " },
    @{ Row = 39; Col = 5; Text = "  -0.07%  " },
    @{ Row = 40; Col = 5; Text = "  -0.29%  " },
    @{ Row = 41; Col = 5; Text = "  +4.69%  " },
    @{ Row = 42; Col = 4; Text = "1.542.19" },
    @{ Row = 42; Col = 5; Text = "  +0.50%  " },
    @{ Row = 43; Col = 5; Text = "  +3.40%  " },
    @{ Row = 44; Col = 5; Text = "  -0.78%  " },
    @{ Row = 45; Col = 5; Text = "  +3.37%  " },
    @{ Row = 46; Col = 5; Text = "  +8.38%  " },
    @{ Row = 47; Col = 5; Text = "  +1.17%  " },
    @{ Row = 48; Col = 5; Text = "  -0.64%  " },
    @{ Row = 49; Col = 5; Text = "  +3.09%  " },
    @{ Row = 50; Col = 5; Text = "  +0.93%  " },
    @{ Row = 51; Col = 4; Text = "2.291.59" },
    @{ Row = 51; Col = 5; Text = "  +2.39%  " }
)

foreach ($u in $textUpdates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Text
}

# Price cells whose new text parses as a plain number (e.g. "228.63")
# -- force Text via NumberFormat "@" so Excel keeps the literal digits
# (incl. trailing zeros) instead of coercing to a Double, then restore
# the default "Normal" style so no visible formatting changes.
$numericLookingUpdates = @(
    @{ Row = 5; Col = 4; Text = "228.63" },
    @{ Row = 6; Col = 4; Text = "0.616" },
    @{ Row = 7; Col = 4; Text = "60.37" },
    @{ Row = 10; Col = 4; Text = "0.0840" },
    @{ Row = 14; Col = 4; Text = "21.99" },
    @{ Row = 15; Col = 4; Text = "0.798" },
    @{ Row = 19; Col = 4; Text = "71.60" },
    @{ Row = 22; Col = 4; Text = "227.43" },
    @{ Row = 25; Col = 4; Text = "2.34" },
    @{ Row = 26; Col = 4; Text = "170.91" },
    @{ Row = 28; Col = 4; Text = "0.140" },
    @{ Row = 31; Col = 4; Text = "0.120" },
    @{ Row = 32; Col = 4; Text = "2.37" },
    @{ Row = 34; Col = 4; Text = "4.70" },
    @{ Row = 36; Col = 4; Text = "6.48" },
    @{ Row = 38; Col = 4; Text = "3.61" },
    @{ Row = 40; Col = 4; Text = "18.21" },
    @{ Row = 43; Col = 4; Text = "101.09" },
    @{ Row = 47; Col = 4; Text = "1.12" }
)

foreach ($u in $numericLookingUpdates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Text
    $cell.Style = "Normal"
}
